$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.917.20"
$ws.Range("E2").Value = "  +5.90%  "
$ws.Range("D3").Value = "2.628.91"
$ws.Range("E3").Value = "  +7.04%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'185.43"
$ws.Range("E5").Value = "  +11.36%  "
$ws.Range("D6").Value = "'581.03"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  +3.33%  "
$ws.Range("E9").Value = "  +9.77%  "
$ws.Range("D10").Value = "2.628.25"
$ws.Range("E10").Value = "  +7.14%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("E12").Value = "  +5.34%  "
$ws.Range("D13").Value = "'4.68"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "73.834.57"
$ws.Range("E14").Value = "  +5.99%  "
$ws.Range("D15").Value = "3.102.18"
$ws.Range("E15").Value = "  +6.36%  "
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "'26.17"
$ws.Range("E17").Value = "  +8.85%  "
$ws.Range("D18").Value = "2.615.09"
$ws.Range("E18").Value = "  +5.87%  "
$ws.Range("D19").Value = "'9.11"
$ws.Range("E19").Value = "  +27.92%  "
$ws.Range("D20").Value = "'11.80"
$ws.Range("E20").Value = "  +9.14%  "
$ws.Range("D21").Value = "'364.28"
$ws.Range("E21").Value = "  +6.28%  "
$ws.Range("D22").Value = "'2.28"
$ws.Range("E22").Value = "  +13.30%  "
$ws.Range("D23").Value = "'4.05"
$ws.Range("E23").Value = "  +4.22%  "
$ws.Range("D24").Value = "'6.16"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "'69.44"
$ws.Range("E26").Value = "  +4.44%  "
$ws.Range("D27").Value = "'4.07"
$ws.Range("E27").Value = "  +3.73%  "
$ws.Range("D28").Value = "'9.24"
$ws.Range("E28").Value = "  +8.18%  "
$ws.Range("E29").Value = "  +6.21%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").Value = "0.0₃0927"
$ws.Range("E31").Value = "  +9.17%  "
$ws.Range("D32").Value = "'518.73"
$ws.Range("E32").Value = "  +15.68%  "
$ws.Range("D33").Value = "'1.38"
$ws.Range("E33").Value = "  +9.73%  "
$ws.Range("D34").Value = "'7.60"
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("D35").Value = "'1.73"
$ws.Range("E35").Value = "  +6.03%  "
$ws.Range("D36").Value = "'0.997"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "'162.29"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("E38").Value = "  +7.26%  "
$ws.Range("D39").Value = "'18.99"
$ws.Range("E39").Value = "  +4.74%  "
$ws.Range("D40").Value = "'19.28"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'4.87"
$ws.Range("E42").Value = "  +8.78%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.65"
$ws.Range("E43").Value = "  +6.99%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "'0.324"
$ws.Range("E44").Value = "  +5.41%  "
$ws.Range("D45").Value = "'161.90"
$ws.Range("E45").Value = "  +22.64%  "
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "'1.17"
$ws.Range("E46").Value = "  +7.17%  "
$ws.Range("D47").Value = "'2.34"
$ws.Range("E47").Value = "  +9.03%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'38.87"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("D49").Value = "'0.0840"
$ws.Range("E49").Value = "  +16.23%  "
$ws.Range("D50").Value = "'3.58"
$ws.Range("E50").Value = "  +5.45%  "
$ws.Range("D51").Value = "'0.521"
$ws.Range("E51").Value = "  +6.06%  "
